$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data values between row 2 and row 3 for columns:
# D (Fecha), I (Calidad), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion), P (Precio $/Kg), Q (Kg o Unidades)

$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2

    $cellRow2.Value2 = $val3
    $cellRow3.Value2 = $val2
}
